$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing row (719) into the new rows (720:750).
# Row 719 has styled cells A-L and N, but not M (the quiz has two mutually
# exclusive follow-up columns, M and N, depending on an earlier answer).
$ws.Range("A719:N719").Copy()
$ws.Range("A720:N750").PasteSpecial(-4122)

# Row 720
$ws.Range("A720").Value = 45200.671034479165
$ws.Range("B720").Value = "wlalsdlcjm4@naver.com"
$ws.Range("C720").Value = "간호학과"
$ws.Range("D720").Value = 20236296
$ws.Range("E720").Value = "차지민"
$ws.Range("F720").Value = "실제로 현장에 나가서 수확량을 파악하고 등급을 매기는 답험(踏驗)을 하였다."
$ws.Range("G720").Value = 0.5
$ws.Range("H720").Value = "6:4"
$ws.Range("I720").Value = "15분의 1"
$ws.Range("J720").Value = "44만호, 153만명"
$ws.Range("K720").Value = "경상"
$ws.Range("L720").Value = "Black"
$ws.Range("N720").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."
$ws.Range("M720").Clear()

# Row 721
$ws.Range("A721").Value = 45200.676602175925
$ws.Range("B721").Value = "qwer030815@naver.com"
$ws.Range("C721").Value = "바이오메디컬"
$ws.Range("D721").Value = 20223639
$ws.Range("E721").Value = "정은진"
$ws.Range("F721").Value = "실제로 현장에 나가서 수확량을 파악하고 등급을 매기는 답험(踏驗)을 하였다."
$ws.Range("G721").Value = 0.1
$ws.Range("H721").Value = "6:4"
$ws.Range("I721").Value = "10분의 1"
$ws.Range("J721").Value = "20만호, 69만명"
$ws.Range("K721").Value = "평안"
$ws.Range("L721").Value = "Black"
$ws.Range("N721").Value = "모름/무응답"
$ws.Range("M721").Clear()

# Row 722
$ws.Range("L722").Copy()
$ws.Range("M722").PasteSpecial(-4122)
$ws.Range("A722").Value = 45200.67752693287
$ws.Range("B722").Value = "duddms5818@naver.com"
$ws.Range("C722").Value = "광고홍보학과"
$ws.Range("D722").Value = 20192638
$ws.Range("E722").Value = "한영은"
$ws.Range("F722").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G722").Value = 0.1
$ws.Range("H722").Value = "6:4"
$ws.Range("I722").Value = "20분의 1"
$ws.Range("J722").Value = "20만호, 69만명"
$ws.Range("K722").Value = "충청"
$ws.Range("L722").Value = "Red"
$ws.Range("M722").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."
$ws.Range("N722").Clear()

# Row 723
$ws.Range("L723").Copy()
$ws.Range("M723").PasteSpecial(-4122)
$ws.Range("A723").Value = 45200.67929005787
$ws.Range("B723").Value = "rudqh0501@gmail.com"
$ws.Range("C723").Value = "화학과"
$ws.Range("D723").Value = 20233417
$ws.Range("E723").Value = "이경보"
$ws.Range("F723").Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Range("G723").Value = 0.5
$ws.Range("H723").Value = "5:5"
$ws.Range("I723").Value = "15분의 1"
$ws.Range("J723").Value = "20만호, 69만명"
$ws.Range("K723").Value = "전라"
$ws.Range("L723").Value = "Red"
$ws.Range("M723").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."
$ws.Range("N723").Clear()

# Row 724
$ws.Range("L724").Copy()
$ws.Range("M724").PasteSpecial(-4122)
$ws.Range("A724").Value = 45200.68185594908
$ws.Range("B724").Value = "lma240228@gmail.com"
$ws.Range("C724").Value = "미래융합스쿨"
$ws.Range("D724").Value = 20236608
$ws.Range("E724").Value = "김재민"
$ws.Range("F724").Value = "‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다."
$ws.Range("G724").Value = 0.3
$ws.Range("H724").Value = "4:6"
$ws.Range("I724").Value = "20분의 1"
$ws.Range("J724").Value = "44만호, 153만명"
$ws.Range("K724").Value = "경기"
$ws.Range("L724").Value = "Red"
$ws.Range("M724").Value = "모름/무응답"
$ws.Range("N724").Clear()

# Row 725
$ws.Range("A725").Value = 45200.682556678235
$ws.Range("B725").Value = "dorud030406@naver.com"
$ws.Range("C725").Value = "미디어스쿨"
$ws.Range("D725").Value = 20222559
$ws.Range("E725").Value = "이예경"
$ws.Range("F725").Value = "과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."
$ws.Range("G725").Value = 0.3
$ws.Range("H725").Value = "6:4"
$ws.Range("I725").Value = "15분의 1"
$ws.Range("J725").Value = "20만호, 69만명"
$ws.Range("K725").Value = "경상"
$ws.Range("L725").Value = "Black"
$ws.Range("N725").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."
$ws.Range("M725").Clear()

# Row 726
$ws.Range("L726").Copy()
$ws.Range("M726").PasteSpecial(-4122)
$ws.Range("A726").Value = 45200.68288702546
$ws.Range("B726").Value = "limyoon0725@daum.net"
$ws.Range("C726").Value = "사회복지학과"
$ws.Range("D726").Value = 20222347
$ws.Range("E726").Value = "임윤서"
$ws.Range("F726").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G726").Value = 0.1
$ws.Range("H726").Value = "6:4"
$ws.Range("I726").Value = "20분의 1"
$ws.Range("J726").Value = "20만호, 69만명"
$ws.Range("K726").Value = "충청"
$ws.Range("L726").Value = "Red"
$ws.Range("M726").Value = "반대한다."
$ws.Range("N726").Clear()

# Row 727
$ws.Range("A727").Value = 45200.68478565972
$ws.Range("B727").Value = "gahee021911@gmail.com"
$ws.Range("C727").Value = "경영대학"
$ws.Range("D727").Value = 20233023
$ws.Range("E727").Value = "임가희"
$ws.Range("F727").Value = "실제로 현장에 나가서 수확량을 파악하고 등급을 매기는 답험(踏驗)을 하였다."
$ws.Range("G727").Value = 0.1
$ws.Range("H727").Value = "6:4"
$ws.Range("I727").Value = "10분의 1"
$ws.Range("J727").Value = "20만호, 69만명"
$ws.Range("K727").Value = "평안"
$ws.Range("L727").Value = "Black"
$ws.Range("N727").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."
$ws.Range("M727").Clear()

# Row 728
$ws.Range("L728").Copy()
$ws.Range("M728").PasteSpecial(-4122)
$ws.Range("A728").Value = 45200.68725362269
$ws.Range("B728").Value = "dungunfight9@gmail.com"
$ws.Range("C728").Value = "금융재무학과"
$ws.Range("D728").Value = 20222986
$ws.Range("E728").Value = "유홍현"
$ws.Range("F728").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G728").Value = 0.3
$ws.Range("H728").Value = "6:4"
$ws.Range("I728").Value = "20분의 1"
$ws.Range("J728").Value = "20만호, 69만명"
$ws.Range("K728").Value = "충청"
$ws.Range("L728").Value = "Red"
$ws.Range("M728").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."
$ws.Range("N728").Clear()

# Row 729
$ws.Range("L729").Copy()
$ws.Range("M729").PasteSpecial(-4122)
$ws.Range("A729").Value = 45200.68982090278
$ws.Range("B729").Value = "tjwls7140@naver.com"
$ws.Range("C729").Value = "반도체디스플레이스쿨"
$ws.Range("D729").Value = 20233353
$ws.Range("E729").Value = "황서진"
$ws.Range("F729").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G729").Value = 0.1
$ws.Range("H729").Value = "6:4"
$ws.Range("I729").Value = "10분의 1"
$ws.Range("J729").Value = "20만호, 69만명"
$ws.Range("K729").Value = "충청"
$ws.Range("L729").Value = "Red"
$ws.Range("M729").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."
$ws.Range("N729").Clear()

# Row 730
$ws.Range("L730").Copy()
$ws.Range("M730").PasteSpecial(-4122)
$ws.Range("A730").Value = 45200.691276516205
$ws.Range("B730").Value = "h20221203@hallym.glab.ac.kr"
$ws.Range("C730").Value = "영어영문학과"
$ws.Range("D730").Value = 20221203
$ws.Range("E730").Value = "권민주"
$ws.Range("F730").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G730").Value = 0.5
$ws.Range("H730").Value = "6:4"
$ws.Range("I730").Value = "10분의 1"
$ws.Range("J730").Value = "20만호, 69만명"
$ws.Range("K730").Value = "경상"
$ws.Range("L730").Value = "Red"
$ws.Range("M730").Value = "모름/무응답"
$ws.Range("N730").Clear()

# Row 731
$ws.Range("L731").Copy()
$ws.Range("M731").PasteSpecial(-4122)
$ws.Range("A731").Value = 45200.692609895836
$ws.Range("B731").Value = "hyj4213@naver.com"
$ws.Range("C731").Value = "미디어스쿨"
$ws.Range("D731").Value = 20232590
$ws.Range("E731").Value = "함영준"
$ws.Range("F731").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G731").Value = 0.9
$ws.Range("H731").Value = "6:4"
$ws.Range("I731").Value = "20분의 1"
$ws.Range("J731").Value = "20만호, 69만명"
$ws.Range("K731").Value = "충청"
$ws.Range("L731").Value = "Red"
$ws.Range("M731").Value = "반대한다."
$ws.Range("N731").Clear()

# Row 732
$ws.Range("L732").Copy()
$ws.Range("M732").PasteSpecial(-4122)
$ws.Range("A732").Value = 45200.69353521991
$ws.Range("B732").Value = "jeh9599@gmail.com"
$ws.Range("C732").Value = "소프트웨어"
$ws.Range("D732").Value = 20235256
$ws.Range("E732").Value = "정은혁"
$ws.Range("F732").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G732").Value = 0.7
$ws.Range("H732").Value = "6:4"
$ws.Range("I732").Value = "20분의 1"
$ws.Range("J732").Value = "44만호, 153만명"
$ws.Range("K732").Value = "평안"
$ws.Range("L732").Value = "Red"
$ws.Range("M732").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."
$ws.Range("N732").Clear()

# Row 733
$ws.Range("L733").Copy()
$ws.Range("M733").PasteSpecial(-4122)
$ws.Range("A733").Value = 45200.69612866898
$ws.Range("B733").Value = "wasabi0724@naver.com"
$ws.Range("C733").Value = "미래융합스쿨"
$ws.Range("D733").Value = 20236645
$ws.Range("E733").Value = "홍지윤"
$ws.Range("F733").Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Range("G733").Value = 0.3
$ws.Range("H733").Value = "4:6"
$ws.Range("I733").Value = "10분의 1"
$ws.Range("J733").Value = "44만호, 153만명"
$ws.Range("K733").Value = "평안"
$ws.Range("L733").Value = "Red"
$ws.Range("M733").Value = "모름/무응답"
$ws.Range("N733").Clear()

# Row 734
$ws.Range("A734").Value = 45200.69644321759
$ws.Range("B734").Value = "chdcks7359@gmail.com"
$ws.Range("C734").Value = "경영학과"
$ws.Range("D734").Value = 20233027
$ws.Range("E734").Value = "장총찬"
$ws.Range("F734").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G734").Value = 0.5
$ws.Range("H734").Value = "4:6"
$ws.Range("I734").Value = "20분의 1"
$ws.Range("J734").Value = "44만호, 153만명"
$ws.Range("K734").Value = "전라"
$ws.Range("L734").Value = "Black"
$ws.Range("N734").Value = "찬성한다."
$ws.Range("M734").Clear()

# Row 735
$ws.Range("A735").Value = 45200.69692814814
$ws.Range("B735").Value = "pungnam04@gmail.com"
$ws.Range("C735").Value = "소프트웨어학부"
$ws.Range("D735").Value = 20235242
$ws.Range("E735").Value = "이종민"
$ws.Range("F735").Value = "과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."
$ws.Range("G735").Value = 0.5
$ws.Range("H735").Value = "6:4"
$ws.Range("I735").Value = "20분의 1"
$ws.Range("J735").Value = "44만호, 153만명"
$ws.Range("K735").Value = "경상"
$ws.Range("L735").Value = "Black"
$ws.Range("N735").Value = "모름/무응답"
$ws.Range("M735").Clear()

# Row 736
$ws.Range("A736").Value = 45200.71004925926
$ws.Range("B736").Value = "wjdcofla24@naver.com"
$ws.Range("C736").Value = "러시아학과"
$ws.Range("D736").Value = 20211726
$ws.Range("E736").Value = "정채림"
$ws.Range("F736").Value = "과전법 체제에서 전국 토지를 세 등급으로 나누고 실제 수확량을 확인하여 징수하였다."
$ws.Range("G736").Value = 0.7
$ws.Range("H736").Value = "6:4"
$ws.Range("I736").Value = "15분의 1"
$ws.Range("J736").Value = "20만호, 69만명"
$ws.Range("K736").Value = "평안"
$ws.Range("L736").Value = "Black"
$ws.Range("N736").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."
$ws.Range("M736").Clear()

# Row 737
$ws.Range("L737").Copy()
$ws.Range("M737").PasteSpecial(-4122)
$ws.Range("A737").Value = 45200.715234675925
$ws.Range("B737").Value = "dosilver1107@naver.com"
$ws.Range("C737").Value = "사회학과"
$ws.Range("D737").Value = 20212233
$ws.Range("E737").Value = "임도은"
$ws.Range("F737").Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Range("G737").Value = 0.7
$ws.Range("H737").Value = "4:6"
$ws.Range("I737").Value = "15분의 1"
$ws.Range("J737").Value = "44만호, 153만명"
$ws.Range("K737").Value = "경상"
$ws.Range("L737").Value = "Red"
$ws.Range("M737").Value = "근로시간과 휴무를 유연하게 조정할 수 있어 찬성한다."
$ws.Range("N737").Clear()

# Row 738
$ws.Range("A738").Value = 45200.715887847226
$ws.Range("B738").Value = "qwerty052799@gmail.com"
$ws.Range("C738").Value = "인공지능융합학부"
$ws.Range("D738").Value = 20236722
$ws.Range("E738").Value = "박성현"
$ws.Range("F738").Value = "‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다."
$ws.Range("G738").Value = 0.9
$ws.Range("H738").Value = "3:7"
$ws.Range("I738").Value = "10분의 1"
$ws.Range("J738").Value = "44만호, 153만명"
$ws.Range("K738").Value = "경상"
$ws.Range("L738").Value = "Black"
$ws.Range("N738").Value = "모름/무응답"
$ws.Range("M738").Clear()

# Row 739
$ws.Range("L739").Copy()
$ws.Range("M739").PasteSpecial(-4122)
$ws.Range("A739").Value = 45200.717868645836
$ws.Range("B739").Value = "tjdbs6201305@naver.com"
$ws.Range("C739").Value = "미디어스쿨"
$ws.Range("D739").Value = 20232545
$ws.Range("E739").Value = "송서윤"
$ws.Range("F739").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G739").Value = 0.3
$ws.Range("H739").Value = "6:4"
$ws.Range("I739").Value = "10분의 1"
$ws.Range("J739").Value = "20만호, 69만명"
$ws.Range("K739").Value = "전라"
$ws.Range("L739").Value = "Red"
$ws.Range("M739").Value = "반대한다."
$ws.Range("N739").Clear()

# Row 740
$ws.Range("A740").Value = 45200.71845825232
$ws.Range("B740").Value = "dbwjdrms21@naver.com"
$ws.Range("C740").Value = "식품영양학과"
$ws.Range("D740").Value = 20183829
$ws.Range("E740").Value = "유정근"
$ws.Range("F740").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G740").Value = 0.1
$ws.Range("H740").Value = "6:4"
$ws.Range("I740").Value = "20분의 1"
$ws.Range("J740").Value = "20만호, 69만명"
$ws.Range("K740").Value = "충청"
$ws.Range("L740").Value = "Black"
$ws.Range("N740").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."
$ws.Range("M740").Clear()

# Row 741
$ws.Range("A741").Value = 45200.72100321759
$ws.Range("B741").Value = "rudtjraudwls@naver.com"
$ws.Range("C741").Value = "간호학과"
$ws.Range("D741").Value = 20236302
$ws.Range("E741").Value = "한경석"
$ws.Range("F741").Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Range("G741").Value = 0.7
$ws.Range("H741").Value = "6:4"
$ws.Range("I741").Value = "20분의 1"
$ws.Range("J741").Value = "44만호, 153만명"
$ws.Range("K741").Value = "평안"
$ws.Range("L741").Value = "Black"
$ws.Range("N741").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."
$ws.Range("M741").Clear()

# Row 742
$ws.Range("A742").Value = 45200.72329980324
$ws.Range("B742").Value = "plzmxn@naver.com"
$ws.Range("C742").Value = "광고홍보학과"
$ws.Range("D742").Value = 20232642
$ws.Range("E742").Value = "현상희"
$ws.Range("F742").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G742").Value = 0.1
$ws.Range("H742").Value = "6:4"
$ws.Range("I742").Value = "10분의 1"
$ws.Range("J742").Value = "20만호, 69만명"
$ws.Range("K742").Value = "평안"
$ws.Range("L742").Value = "Black"
$ws.Range("N742").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."
$ws.Range("M742").Clear()

# Row 743
$ws.Range("L743").Copy()
$ws.Range("M743").PasteSpecial(-4122)
$ws.Range("A743").Value = 45200.72627798611
$ws.Range("B743").Value = "sdw0820@gmail.com"
$ws.Range("C743").Value = "식품영양학과"
$ws.Range("D743").Value = 20233819
$ws.Range("E743").Value = "성동휘"
$ws.Range("F743").Value = "‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다."
$ws.Range("G743").Value = 0.1
$ws.Range("H743").Value = "6:4"
$ws.Range("I743").Value = "20분의 1"
$ws.Range("J743").Value = "44만호, 153만명"
$ws.Range("K743").Value = "충청"
$ws.Range("L743").Value = "Red"
$ws.Range("M743").Value = "모름/무응답"
$ws.Range("N743").Clear()

# Row 744
$ws.Range("A744").Value = 45200.72919329861
$ws.Range("B744").Value = "041030top@naver.com"
$ws.Range("C744").Value = "데이터사이언스"
$ws.Range("D744").Value = 20233257
$ws.Range("E744").Value = "최영국"
$ws.Range("F744").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G744").Value = 0.1
$ws.Range("H744").Value = "6:4"
$ws.Range("I744").Value = "20분의 1"
$ws.Range("J744").Value = "20만호, 69만명"
$ws.Range("K744").Value = "충청"
$ws.Range("L744").Value = "Black"
$ws.Range("N744").Value = "찬성한다."
$ws.Range("M744").Clear()

# Row 745
$ws.Range("A745").Value = 45200.732332534724
$ws.Range("B745").Value = "codmsrjf@naver.com"
$ws.Range("C745").Value = "콘텐츠 IT"
$ws.Range("D745").Value = 20205253
$ws.Range("E745").Value = "정채은"
$ws.Range("F745").Value = "‘세(稅)’는 사전의 소유자가 국가에 상납하는 지대를 뜻한다."
$ws.Range("G745").Value = 0.3
$ws.Range("H745").Value = "6:4"
$ws.Range("I745").Value = "20분의 1"
$ws.Range("J745").Value = "20만호, 69만명"
$ws.Range("K745").Value = "평안"
$ws.Range("L745").Value = "Black"
$ws.Range("N745").Value = "찬성한다."
$ws.Range("M745").Clear()

# Row 746
$ws.Range("L746").Copy()
$ws.Range("M746").PasteSpecial(-4122)
$ws.Range("A746").Value = 45200.732799375
$ws.Range("B746").Value = "20232327@hallym.ac.kr"
$ws.Range("C746").Value = "사회복지학부 "
$ws.Range("D746").Value = 20232327
$ws.Range("E746").Value = "박소희 "
$ws.Range("F746").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G746").Value = 0.1
$ws.Range("H746").Value = "6:4"
$ws.Range("I746").Value = "20분의 1"
$ws.Range("J746").Value = "20만호, 69만명"
$ws.Range("K746").Value = "충청"
$ws.Range("L746").Value = "Red"
$ws.Range("M746").Value = "모름/무응답"
$ws.Range("N746").Clear()

# Row 747
$ws.Range("A747").Value = 45200.736784733796
$ws.Range("B747").Value = "sungyeon0803@gmail.com"
$ws.Range("C747").Value = "정치행정학과"
$ws.Range("D747").Value = 20232437
$ws.Range("E747").Value = "최성연"
$ws.Range("F747").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G747").Value = 0.1
$ws.Range("H747").Value = "6:4"
$ws.Range("I747").Value = "15분의 1"
$ws.Range("J747").Value = "44만호, 153만명"
$ws.Range("K747").Value = "평안"
$ws.Range("L747").Value = "Black"
$ws.Range("N747").Value = "노동자가 과도한 연장근로를 받을 수 있어 반대한다."
$ws.Range("M747").Clear()

# Row 748
$ws.Range("L748").Copy()
$ws.Range("M748").PasteSpecial(-4122)
$ws.Range("A748").Value = 45200.751256550924
$ws.Range("B748").Value = "cba05049@naver.com"
$ws.Range("C748").Value = "법학과"
$ws.Range("D748").Value = 20172720
$ws.Range("E748").Value = "박상우"
$ws.Range("F748").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G748").Value = 0.1
$ws.Range("H748").Value = "6:4"
$ws.Range("I748").Value = "20분의 1"
$ws.Range("J748").Value = "20만호, 69만명"
$ws.Range("K748").Value = "충청"
$ws.Range("L748").Value = "Red"
$ws.Range("M748").Value = "반대한다."
$ws.Range("N748").Clear()

# Row 749
$ws.Range("L749").Copy()
$ws.Range("M749").PasteSpecial(-4122)
$ws.Range("A749").Value = 45200.76741747685
$ws.Range("B749").Value = "aktnftk3520@naver.com"
$ws.Range("C749").Value = "미래융합스쿨"
$ws.Range("D749").Value = 20236628
$ws.Range("E749").Value = "어수련"
$ws.Range("F749").Value = "‘조(租)’는 공전(公田)의 경작자가 국고에 상납하는 지대 또는 사전(私田)의 경작자가 전주에게 바치는 지대를 뜻한다."
$ws.Range("G749").Value = 0.3
$ws.Range("H749").Value = "7:3"
$ws.Range("I749").Value = "10분의 1"
$ws.Range("J749").Value = "15만호,  32만명"
$ws.Range("K749").Value = "충청"
$ws.Range("L749").Value = "Red"
$ws.Range("M749").Value = "반대한다."
$ws.Range("N749").Clear()

# Row 750
$ws.Range("A750").Value = 45200.76876439815
$ws.Range("B750").Value = "dlwlgy0001@naver.com"
$ws.Range("C750").Value = "경영학과"
$ws.Range("D750").Value = 20203024
$ws.Range("E750").Value = "이지효"
$ws.Range("F750").Value = "등급에 따라 일정한 비율로 세금을 감면해 주는 대동(大同)법을 실시하였다."
$ws.Range("G750").Value = 0.1
$ws.Range("H750").Value = "6:4"
$ws.Range("I750").Value = "20분의 1"
$ws.Range("J750").Value = "20만호, 69만명"
$ws.Range("K750").Value = "충청"
$ws.Range("L750").Value = "Black"
$ws.Range("N750").Value = "찬성한다."
$ws.Range("M750").Clear()
